# Indian MF 1st Stab - add 9 new weekly columns (Jun_16 .. Sep_08) to the
# MarketBeat rank tracker sheet, inserting them right after column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before the old column B; this shifts all the existing
# weekly data (previously B:R) to the right by 9 columns (now K:AA), exactly
# like Excel's normal "Insert Copied/Blank Columns" behaviour.
$ws.Columns("B:J").Insert()

# New week-ending date headers for the newly inserted columns (most recent
# week first, mirroring the existing right-to-left chronological layout).
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# Set column widths for the 9 new columns to match the rest of the data
# columns (27.0), same as every other weekly column on the sheet.
$ws.Range("B1:J1").EntireColumn.ColumnWidth = 27

# Fill the new weekly cells for every analyst row with "UN" (unchanged /
# no new rating that week), same placeholder used throughout the sheet.
For ($r = 2; $r -le 33; $r++) {
    For ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}

# BidaskClub (row 22) actually published a new rating during the Jun_16
# week: a downgrade from Hold to Sell on 6/13/2019.
$ws.Range("J22").Value = "6/13/2019,Downgrades,Hold -> Sell,"
